$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction1")

$values = @(0, 2, 3, 4, 5, 6, 7, 8, 10, 12, 13, 14, 15, 16, 19, 20, 21, 22, 24, 25, 27, 28, 29, 30, 31, 32, 33, 34, 38, 39, 41, 42, 43, 44, 46, 47, 51, 54, 55, 57, 58, 60, 62, 63, 64, 65, 67, 68, 69, 70, 71)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
